$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.550.76"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "3.812.92"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "684.31"
$ws.Range("E5").Value = "  +9.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.45"
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("D7").Value = "3.812.03"
$ws.Range("E7").Value = "  +1.38%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  +7.56%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.05"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "4.451.76"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "3.806.81"
$ws.Range("E16").Value = "  +1.03%  "
$ws.Range("D17").Value = "70.602.47"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.72"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("E19").Value = "  +2.51%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.39"
$ws.Range("E21").Value = "  +20.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "475.77"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.715"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.56"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.29"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "3.960.82"
$ws.Range("E30").Value = "  +1.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.95"
$ws.Range("E31").Value = "  +10.79%  "
$ws.Range("E32").Value = "  +3.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.41"
$ws.Range("E33").Value = "  +4.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.71"
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("E35").Value = "  +6.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.15"
$ws.Range("E36").Value = "  +2.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "3.759.37"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.966"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +13.67%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "46.37"
$ws.Range("E46").Value = "  +8.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "159.88"
$ws.Range("E47").Value = "  +2.99%  "
$ws.Range("E48").Value = "  +7.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000294"
$ws.Range("E50").Value = "  +9.59%  "
$ws.Range("E51").Value = "  +2.13%  "
